# Add a new "Dist" column (H) to Sheet1 classifying each row as WT (Col-0)
# or Mutant (all other genotypes), and update the sheet view (zoom/scroll/
# selection) to match the reviewer's last-saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("H1").Value = "Dist"

# Rows 2-13 are the Col-0 (wild type) samples
$ws.Range("H2:H13").Value = "WT"

# Rows 14-85 are the various mutant genotypes (dwf4, cpd, det2, brox1,2, bri116)
$ws.Range("H14:H85").Value = "Mutant"

# Match the saved view state: scrolled to row 46, zoomed to 85%, with
# H14:H85 selected (active cell H14)
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H14:H85").Select() | Out-Null
